$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.793.24"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "3.532.51"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'604.83"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'196.29"
$ws.Range("E6").Value = "  +5.58%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "'0.650"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'53.67"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "'0.0000304"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "'9.50"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "4.090.24"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "'603.47"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "69.994.73"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "'19.09"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "'12.75"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "3.537.29"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").Value = "'0.994"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'18.34"
$ws.Range("E22").Value = "  +5.47%  "
$ws.Range("D23").Value = "'5.29"
$ws.Range("E23").Value = "  +5.00%  "
$ws.Range("D24").Value = "'102.23"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").Value = "'4.63"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  +3.26%  "
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.12"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").Value = "'4.30"
$ws.Range("E31").Value = "  +13.78%  "
$ws.Range("D32").Value = "'12.51"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "'63.18"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("D35").Value = "0.0₃0878"
$ws.Range("E35").Value = "  +14.53%  "
$ws.Range("D36").Value = "3.726.45"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'3.05"
$ws.Range("E38").Value = "  -2.91%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").Value = "'36.62"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "'487.05"
$ws.Range("E42").Value = "  -7.39%  "
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "'0.0456"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").Value = "'8.56"
$ws.Range("E49").Value = "  -3.42%  "
$ws.Range("D50").Value = "'0.000255"
$ws.Range("E50").Value = "  +5.99%  "
$ws.Range("D51").Value = "'131.17"
$ws.Range("E51").Value = "  -0.32%  "
